# Add a "Team" column (J) to the roster sheet:
#   - J1 header "Team", bold, centered/top-aligned, with left+right thin borders
#   - J2:J92 filled with the constant value "Lib" for every player row
# (commit: "added foreign key + check")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row on the sheet (header + 91 player rows -> row 92)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# --- Header cell J1 -----------------------------------------------------
$header = $ws.Range("J1")
$header.Value = "Team"

$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop

# Left/right thin borders only (matches the style used for the other
# header cells, but without top/bottom rules)
$header.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$header.Borders.Item(7).LineStyle = 1    # xlEdgeLeft

# --- Data cells J2:J92 ---------------------------------------------------
$dataRange = $ws.Range("J2:J" + $lastRow)
$dataRange.Value = "Lib"

# Leave the new column's data selected, mirroring the state the workbook
# was saved in after the edit.
$ws.Range("A2:J" + $lastRow).Select() | Out-Null

Write-Output "added Team column (header + $($lastRow - 1) data rows)"
